$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at 43-44, shifting existing rows 43-118 down by 2
$ws.Rows("43:44").Insert()

# Row 43
$ws.Range("A43").Value = 9
$ws.Range("B43").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C43").Value = 'Metropolitana'
$ws.Range("D43").Value = 44495
$ws.Range("E43").Value = 13
$ws.Range("F43").Value = 100112003
$ws.Range("G43").Value = 'Ajo'
$ws.Range("H43").Value = 'Rosado'
$ws.Range("I43").Value = '1a nueva(o)'
$ws.Range("J43").Value = 130
$ws.Range("K43").Value = 3500
$ws.Range("L43").Value = 3500
$ws.Range("M43").Value = 3500
$ws.Range("N43").Value = '$/paquete 20 unidades (volumen en unidades)'
$ws.Range("O43").Value = 'Provincia de Talagante'
$ws.Range("P43").Value = 175
$ws.Range("Q43").Value = 20
$ws.Range("R43").Value = 'Hortaliza'

# Row 44
$ws.Range("A44").Value = 9
$ws.Range("B44").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C44").Value = 'Metropolitana'
$ws.Range("D44").Value = 44495
$ws.Range("E44").Value = 13
$ws.Range("F44").Value = 100112003
$ws.Range("G44").Value = 'Ajo'
$ws.Range("H44").Value = 'Rosado'
$ws.Range("I44").Value = 'Extra nueva (o)'
$ws.Range("J44").Value = 60
$ws.Range("K44").Value = 4000
$ws.Range("L44").Value = 4000
$ws.Range("M44").Value = 4000
$ws.Range("N44").Value = '$/paquete 20 unidades (volumen en unidades)'
$ws.Range("O44").Value = 'Provincia de Talagante'
$ws.Range("P44").Value = 200
$ws.Range("Q44").Value = 20
$ws.Range("R44").Value = 'Hortaliza'
